$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tâches")
$ws.Activate()

# Row 10: Librairie de lecture du capteur de luminosité -> Alexandre Bodart, En cours
$ws.Range("F10").Formula = "=Participants!A3"
$ws.Range("G10").Value = "En cours"

# Row 11: Librairie de lecture du capteur de son -> Alexandre Bodart, En cours
$ws.Range("F11").Formula = "=Participants!A3"
$ws.Range("G11").Value = "En cours"

# Row 12: Librairie de lecture du capteur SensiBLE -> François Monteil, En cours
$ws.Range("F12").Formula = "=Participants!A2"
$ws.Range("G12").Value = "En cours"

# Row 14: Créer le programme -> Les deux, En attente
$ws.Range("F14").Value = "Les deux"
$ws.Range("G14").Value = "En attente"

# Selection / view changes
$ws.Range("C3:H15").Select()
